# Scheduled runner refresh: update market-price derived columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# for the rows whose underlying market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 29414490
$ws.Cells.Item(6, 9).Value = 62500170
$ws.Cells.Item(6, 10).Value = 4998.8887
$ws.Cells.Item(6, 11).Value = 187500510
$ws.Cells.Item(6, 12).Value = 14996.6661
$ws.Cells.Item(6, 13).Value = -187500398
$ws.Cells.Item(6, 14).Value = -15220.6661
$ws.Cells.Item(12, 8).Value = 368
$ws.Cells.Item(12, 10).Value = 440.66666
$ws.Cells.Item(12, 12).Value = 440.66666
$ws.Cells.Item(12, 14).Value = -780.66666
$ws.Cells.Item(19, 8).Value = 1168.7778
$ws.Cells.Item(19, 10).Value = 1603.1666
$ws.Cells.Item(19, 12).Value = 1603.1666
$ws.Cells.Item(19, 14).Value = -1953.1666
$ws.Cells.Item(70, 8).Value = 52386150
$ws.Cells.Item(70, 10).Value = 83340440
$ws.Cells.Item(70, 12).Value = 250021320
$ws.Cells.Item(70, 14).Value = -250021860
$ws.Cells.Item(73, 8).Value = 52386150
$ws.Cells.Item(73, 10).Value = 83340440
$ws.Cells.Item(73, 12).Value = 250021320
$ws.Cells.Item(73, 14).Value = -250023192
$ws.Cells.Item(74, 8).Value = 3828.1667
$ws.Cells.Item(74, 9).Value = 3828.1667
$ws.Cells.Item(74, 11).Value = 3828.1667
$ws.Cells.Item(74, 13).Value = -2892.1667
$ws.Cells.Item(76, 8).Value = 3910
$ws.Cells.Item(76, 10).Value = 4975
$ws.Cells.Item(76, 12).Value = 4975
$ws.Cells.Item(76, 14).Value = -5605
$ws.Cells.Item(77, 8).Value = 3828.1667
$ws.Cells.Item(77, 9).Value = 3828.1667
$ws.Cells.Item(77, 11).Value = 19140.8335
$ws.Cells.Item(77, 13).Value = -14460.8335
$ws.Cells.Item(79, 8).Value = 3910
$ws.Cells.Item(79, 10).Value = 4975
$ws.Cells.Item(79, 12).Value = 4975
$ws.Cells.Item(79, 14).Value = -7159
$ws.Cells.Item(116, 8).Value = 3324.75
$ws.Cells.Item(116, 9).Value = 3500
$ws.Cells.Item(116, 10).Value = 3149.5
$ws.Cells.Item(116, 11).Value = 3500
$ws.Cells.Item(116, 12).Value = 3149.5
$ws.Cells.Item(116, 13).Value = -58
$ws.Cells.Item(116, 14).Value = -10033.5
$ws.Cells.Item(129, 8).Value = 1728.4445
$ws.Cells.Item(129, 9).Value = 914.4
$ws.Cells.Item(129, 10).Value = 2041.5385
$ws.Cells.Item(129, 11).Value = 2743.2
$ws.Cells.Item(129, 12).Value = 6124.6155
$ws.Cells.Item(129, 13).Value = 2256.8
$ws.Cells.Item(129, 14).Value = -16124.6155
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = $null
$ws.Cells.Item(132, 14).Value = $null
$ws.Cells.Item(137, 8).Value = 3381.1667
$ws.Cells.Item(137, 9).Value = 2634
$ws.Cells.Item(137, 10).Value = 4128.3335
$ws.Cells.Item(137, 11).Value = 7902
$ws.Cells.Item(137, 12).Value = 12385.0005
$ws.Cells.Item(137, 13).Value = -5352
$ws.Cells.Item(137, 14).Value = -17485.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1282.5714
$ws.Cells.Item(45, 9).Value = 1188.9231
$ws.Cells.Item(45, 11).Value = 1188.9231
$ws.Cells.Item(45, 13).Value = -811.9231
$ws.Cells.Item(63, 8).Value = 5400.154
$ws.Cells.Item(63, 9).Value = 3355.7778
$ws.Cells.Item(63, 11).Value = 3355.7778
$ws.Cells.Item(63, 13).Value = -2669.7778
$ws.Cells.Item(66, 8).Value = 5400.154
$ws.Cells.Item(66, 9).Value = 3355.7778
$ws.Cells.Item(66, 11).Value = 16778.889
$ws.Cells.Item(66, 13).Value = -13346.889
$ws.Cells.Item(74, 8).Value = 4354.5454
$ws.Cells.Item(74, 9).Value = 4388.6
$ws.Cells.Item(74, 11).Value = 4388.6
$ws.Cells.Item(74, 13).Value = -3514.6
$ws.Cells.Item(77, 8).Value = 4354.5454
$ws.Cells.Item(77, 9).Value = 4388.6
$ws.Cells.Item(77, 11).Value = 21943
$ws.Cells.Item(77, 13).Value = -17575
$ws.Cells.Item(102, 8).Value = 2416
$ws.Cells.Item(102, 9).Value = 1480.7222
$ws.Cells.Item(102, 11).Value = 1480.7222
$ws.Cells.Item(102, 13).Value = 141.2778000000001
$ws.Cells.Item(110, 8).Value = 781.36
$ws.Cells.Item(110, 9).Value = 781.36
$ws.Cells.Item(110, 11).Value = 781.36
$ws.Cells.Item(110, 13).Value = 1263.64
$ws.Cells.Item(122, 8).Value = 6175970.5
$ws.Cells.Item(122, 9).Value = 7939478
$ws.Cells.Item(122, 10).Value = 3693.5
$ws.Cells.Item(122, 11).Value = 23818434
$ws.Cells.Item(122, 12).Value = 11080.5
$ws.Cells.Item(122, 13).Value = -23815984
$ws.Cells.Item(122, 14).Value = -15980.5
$ws.Cells.Item(132, 8).Value = 2013.9231
$ws.Cells.Item(132, 9).Value = 1814.48
$ws.Cells.Item(132, 11).Value = 5443.440000000001
$ws.Cells.Item(132, 13).Value = -2913.440000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1692.9286
$ws.Cells.Item(20, 9).Value = 1880.2
$ws.Cells.Item(20, 11).Value = 1880.2
$ws.Cells.Item(20, 13).Value = -1633.2
$ws.Cells.Item(105, 8).Value = 1790.7273
$ws.Cells.Item(105, 9).Value = 1614.25
$ws.Cells.Item(105, 10).Value = 2261.3333
$ws.Cells.Item(105, 11).Value = 1614.25
$ws.Cells.Item(105, 12).Value = 2261.3333
$ws.Cells.Item(105, 13).Value = 132.75
$ws.Cells.Item(105, 14).Value = -5755.3333
$ws.Cells.Item(107, 8).Value = 8501.223
$ws.Cells.Item(107, 9).Value = 3255.5
$ws.Cells.Item(107, 10).Value = 10000
$ws.Cells.Item(107, 11).Value = 3255.5
$ws.Cells.Item(107, 12).Value = 10000
$ws.Cells.Item(107, 13).Value = -1335.5
$ws.Cells.Item(107, 14).Value = -13840

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 9400
$ws.Cells.Item(62, 9).Value = 11000
$ws.Cells.Item(62, 10).Value = 7000
$ws.Cells.Item(62, 11).Value = 11000
$ws.Cells.Item(62, 12).Value = 7000
$ws.Cells.Item(62, 13).Value = -10376
$ws.Cells.Item(62, 14).Value = -8248
$ws.Cells.Item(65, 8).Value = 9400
$ws.Cells.Item(65, 9).Value = 11000
$ws.Cells.Item(65, 10).Value = 7000
$ws.Cells.Item(65, 11).Value = 55000
$ws.Cells.Item(65, 12).Value = 35000
$ws.Cells.Item(65, 13).Value = -51880
$ws.Cells.Item(65, 14).Value = -41240
$ws.Cells.Item(107, 8).Value = 1452.6875
$ws.Cells.Item(107, 9).Value = 455.14285
$ws.Cells.Item(107, 10).Value = 1732
$ws.Cells.Item(107, 11).Value = 455.14285
$ws.Cells.Item(107, 12).Value = 1732
$ws.Cells.Item(107, 13).Value = 1464.85715
$ws.Cells.Item(107, 14).Value = -5572
$ws.Cells.Item(122, 8).Value = 3773.0588
$ws.Cells.Item(122, 9).Value = 1413.25
$ws.Cells.Item(122, 10).Value = 4499.154
$ws.Cells.Item(122, 11).Value = 4239.75
$ws.Cells.Item(122, 12).Value = 13497.462
$ws.Cells.Item(122, 13).Value = -1789.75
$ws.Cells.Item(122, 14).Value = -18397.462
$ws.Cells.Item(132, 8).Value = 2277.139
$ws.Cells.Item(132, 9).Value = 2275.6
$ws.Cells.Item(132, 10).Value = 2284.8333
$ws.Cells.Item(132, 11).Value = 6826.799999999999
$ws.Cells.Item(132, 12).Value = 6854.499899999999
$ws.Cells.Item(132, 13).Value = -4296.799999999999
$ws.Cells.Item(132, 14).Value = -11914.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(13, 8).Value = 1027.4
$ws.Cells.Item(13, 9).Value = 409.25
$ws.Cells.Item(13, 10).Value = 3500
$ws.Cells.Item(13, 11).Value = 1227.75
$ws.Cells.Item(13, 12).Value = 10500
$ws.Cells.Item(13, 13).Value = -1059.75
$ws.Cells.Item(13, 14).Value = -10836

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 561.375
$ws.Cells.Item(2, 9).Value = 61.363636
$ws.Cells.Item(2, 10).Value = 1661.4
$ws.Cells.Item(2, 11).Value = 61.363636
$ws.Cells.Item(2, 12).Value = 1661.4
$ws.Cells.Item(2, 13).Value = 51.636364
$ws.Cells.Item(2, 14).Value = -1887.4
$ws.Cells.Item(132, 8).Value = 2303.9092
$ws.Cells.Item(132, 9).Value = 2343.7144
$ws.Cells.Item(132, 10).Value = 1468
$ws.Cells.Item(132, 11).Value = 7031.1432
$ws.Cells.Item(132, 12).Value = 4404
$ws.Cells.Item(132, 13).Value = -4501.1432
$ws.Cells.Item(132, 14).Value = -9464

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1334.3914
$ws.Cells.Item(55, 9).Value = 639.0714
$ws.Cells.Item(55, 11).Value = 639.0714
$ws.Cells.Item(55, 13).Value = -466.0714
$ws.Cells.Item(68, 8).Value = 9406.691999999999
$ws.Cells.Item(68, 9).Value = 2499.25
$ws.Cells.Item(68, 10).Value = 12476.667
$ws.Cells.Item(68, 11).Value = 2499.25
$ws.Cells.Item(68, 12).Value = 12476.667
$ws.Cells.Item(68, 13).Value = -1750.25
$ws.Cells.Item(68, 14).Value = -13974.667
$ws.Cells.Item(71, 8).Value = 9406.691999999999
$ws.Cells.Item(71, 9).Value = 2499.25
$ws.Cells.Item(71, 10).Value = 12476.667
$ws.Cells.Item(71, 11).Value = 12496.25
$ws.Cells.Item(71, 12).Value = 62383.335
$ws.Cells.Item(71, 13).Value = -8752.25
$ws.Cells.Item(71, 14).Value = -69871.33499999999
$ws.Cells.Item(100, 8).Value = 10122.375
$ws.Cells.Item(100, 10).Value = 10122.375
$ws.Cells.Item(100, 12).Value = 10122.375
$ws.Cells.Item(100, 14).Value = -11204.375
$ws.Cells.Item(135, 8).Value = 90214.5
$ws.Cells.Item(135, 10).Value = 90214.5
$ws.Cells.Item(135, 12).Value = 90214.5
$ws.Cells.Item(135, 14).Value = -100354.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 10634
$ws.Cells.Item(62, 9).Value = 6500
$ws.Cells.Item(62, 10).Value = 11224.571
$ws.Cells.Item(62, 11).Value = 6500
$ws.Cells.Item(62, 12).Value = 11224.571
$ws.Cells.Item(62, 13).Value = -5876
$ws.Cells.Item(62, 14).Value = -12472.571
$ws.Cells.Item(65, 8).Value = 10634
$ws.Cells.Item(65, 9).Value = 6500
$ws.Cells.Item(65, 10).Value = 11224.571
$ws.Cells.Item(65, 11).Value = 32500
$ws.Cells.Item(65, 12).Value = 56122.855
$ws.Cells.Item(65, 13).Value = -29380
$ws.Cells.Item(65, 14).Value = -62362.855
$ws.Cells.Item(107, 8).Value = 1982.3077
$ws.Cells.Item(107, 9).Value = 2098.8572
$ws.Cells.Item(107, 10).Value = 1846.3334
$ws.Cells.Item(107, 11).Value = 6296.571599999999
$ws.Cells.Item(107, 12).Value = 5539.0002
$ws.Cells.Item(107, 13).Value = -4376.571599999999
$ws.Cells.Item(107, 14).Value = -9379.0002
